# Update "想去人数" (number of people interested) figures in column F
# for the sheets "展览" and "全部类型".
$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 115
    3  = 404
    4  = 11803
    5  = 976
    10 = 174
    13 = 54
    17 = 1411
    19 = 915
    20 = 114
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
